$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2861
$wsExhibit.Range("F5").Value = 17

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 2861
$wsAll.Range("F10").Value = 17
